$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-30, columns A-F
$data = @(
    @(7,42,5,24000,140,60),
    @(7,42,5,24000,140,61),
    @(7,42,5,24000,140,62),
    @(7,42,5,24000,140,63),
    @(7,42,5,24000,140,64),
    @(7,42,5,24000,140,65),
    @(7,42,5,24000,140,50),
    @(14,42,5,24000,140,80),
    @(14,42,5,24000,140,60),
    @(7,32,5,24000,120,80),
    @(7,32,5,24000,120,60),
    @(14,32,5,24000,120,100),
    @(14,32,5,24000,120,80),
    @(7,40,5,28000,140,80),
    @(7,35,5,28000,140,90),
    @(7,30,5,28000,140,120),
    @(7,25,5,28000,140,130),
    @(14,40,5,28000,140,80),
    @(14,35,5,28000,140,90),
    @(14,30,5,28000,140,120),
    @(14,25,5,28000,140,130),
    @(7,40,5,28000,100,60),
    @(7,35,5,28000,100,60),
    @(7,30,5,28000,100,60),
    @(7,25,5,28000,100,60),
    @(14,40,5,20000,90,50),
    @(14,35,5,20000,90,50),
    @(14,30,5,20000,90,50),
    @(14,25,5,20000,90,50)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# Update the active selection to reflect the last edited cell (F26)
$ws.Range("F26").Select()
